$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 748.7143
$ws.Range("I2").Value = 110.5
$ws.Range("K2").Value = 110.5
$ws.Range("M2").Value = 2.5
$ws.Range("H40").Value = 2885.625
$ws.Range("I40").Value = 2024
$ws.Range("J40").Value = 3747.25
$ws.Range("K40").Value = 2024
$ws.Range("L40").Value = 3747.25
$ws.Range("M40").Value = -1849
$ws.Range("N40").Value = -4097.25
$ws.Range("H54").Value = 20076
$ws.Range("I54").Value = 20076
$ws.Range("K54").Value = 20076
$ws.Range("M54").Value = -19590
$ws.Range("H55").Value = 149.6
$ws.Range("I55").Value = 182.66667
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 182.66667
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = 31.33332999999999
$ws.Range("N55").Value = -528
$ws.Range("H92").Value = 1662.6875
$ws.Range("I92").Value = 1551.25
$ws.Range("K92").Value = 1551.25
$ws.Range("M92").Value = -303.25
$ws.Range("H116").Value = 5639.2
$ws.Range("I116").Value = 5125
$ws.Range("J116").Value = 5767.75
$ws.Range("K116").Value = 5125
$ws.Range("L116").Value = 5767.75
$ws.Range("M116").Value = -1683
$ws.Range("N116").Value = -12651.75
$ws.Range("H132").Value = 4493.96
$ws.Range("I132").Value = 3826.1904
$ws.Range("K132").Value = 11478.5712
$ws.Range("M132").Value = -8948.5712
$ws.Range("H137").Value = 2106.8333
$ws.Range("I137").Value = 2168.2
$ws.Range("K137").Value = 6504.599999999999
$ws.Range("M137").Value = -3954.599999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3459.0862
$ws.Range("I32").Value = 762.0769
$ws.Range("K32").Value = 762.0769
$ws.Range("M32").Value = -475.0769
$ws.Range("H61").Value = 4369.3823
$ws.Range("I61").Value = 3592.3704
$ws.Range("J61").Value = 7366.4287
$ws.Range("K61").Value = 3592.3704
$ws.Range("L61").Value = 7366.4287
$ws.Range("M61").Value = -3380.3704
$ws.Range("N61").Value = -7790.4287
$ws.Range("H74").Value = 2441.6155
$ws.Range("I74").Value = 2039.36
$ws.Range("J74").Value = 12498
$ws.Range("K74").Value = 2039.36
$ws.Range("L74").Value = 12498
$ws.Range("M74").Value = -1165.36
$ws.Range("N74").Value = -14246
$ws.Range("H77").Value = 2441.6155
$ws.Range("I77").Value = 2039.36
$ws.Range("J77").Value = 12498
$ws.Range("K77").Value = 10196.8
$ws.Range("L77").Value = 62490
$ws.Range("M77").Value = -5828.799999999999
$ws.Range("N77").Value = -71226
$ws.Range("H132").Value = 2262.652
$ws.Range("I132").Value = 2174.6365
$ws.Range("K132").Value = 6523.9095
$ws.Range("M132").Value = -3993.9095
$ws.Range("H135").Value = 69173.92
$ws.Range("J135").Value = 77522.336
$ws.Range("L135").Value = 77522.336
$ws.Range("N135").Value = -87662.336
$ws.Range("H136").Value = 4369.3823
$ws.Range("I136").Value = 3592.3704
$ws.Range("J136").Value = 7366.4287
$ws.Range("K136").Value = 10777.1112
$ws.Range("L136").Value = 22099.2861
$ws.Range("M136").Value = -8227.111199999999
$ws.Range("N136").Value = -27199.2861

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2516.8867
$ws.Range("I134").Value = 2499.92
$ws.Range("K134").Value = 7499.76
$ws.Range("M134").Value = -4964.76
$ws.Range("H135").Value = 91599.39999999999
$ws.Range("J135").Value = 91599.39999999999
$ws.Range("L135").Value = 91599.39999999999
$ws.Range("N135").Value = -101739.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8211
$ws.Range("I31").Value = 8824.130999999999
$ws.Range("J31").Value = 6800.8
$ws.Range("K31").Value = 8824.130999999999
$ws.Range("L31").Value = 6800.8
$ws.Range("M31").Value = -8529.130999999999
$ws.Range("N31").Value = -7390.8
$ws.Range("H34").Value = 8211
$ws.Range("I34").Value = 8824.130999999999
$ws.Range("J34").Value = 6800.8
$ws.Range("K34").Value = 8824.130999999999
$ws.Range("L34").Value = 6800.8
$ws.Range("M34").Value = -8622.130999999999
$ws.Range("N34").Value = -7204.8
$ws.Range("H58").Value = 3726
$ws.Range("I58").Value = 4551.5
$ws.Range("K58").Value = 4551.5
$ws.Range("M58").Value = -4348.5
$ws.Range("H64").Value = 86666.664
$ws.Range("J64").Value = 86666.664
$ws.Range("L64").Value = 86666.664
$ws.Range("N64").Value = -87162.664
$ws.Range("H67").Value = 86666.664
$ws.Range("J67").Value = 86666.664
$ws.Range("L67").Value = 86666.664
$ws.Range("N67").Value = -88382.664
$ws.Range("H87").Value = 33750
$ws.Range("I87").Value = 17500
$ws.Range("J87").Value = 50000
$ws.Range("K87").Value = 17500
$ws.Range("L87").Value = 50000
$ws.Range("M87").Value = -16314
$ws.Range("N87").Value = -52372
$ws.Range("H88").Value = 31855.5
$ws.Range("J88").Value = 31933.285
$ws.Range("L88").Value = 31933.285
$ws.Range("N88").Value = -32745.285
$ws.Range("H90").Value = 33750
$ws.Range("I90").Value = 17500
$ws.Range("J90").Value = 50000
$ws.Range("K90").Value = 52500
$ws.Range("L90").Value = 150000
$ws.Range("M90").Value = -46572
$ws.Range("N90").Value = -161856
$ws.Range("H91").Value = 31855.5
$ws.Range("J91").Value = 31933.285
$ws.Range("L91").Value = 31933.285
$ws.Range("N91").Value = -34741.285
$ws.Range("H132").Value = 2013.091
$ws.Range("I132").Value = 1793.8889
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 5381.6667
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -2851.6667
$ws.Range("N132").Value = -14058.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 1328
$ws.Range("I134").Value = 988.5
$ws.Range("J134").Value = 2007
$ws.Range("K134").Value = 2965.5
$ws.Range("L134").Value = 6021
$ws.Range("M134").Value = -430.5
$ws.Range("N134").Value = -11091
$ws.Range("H136").Value = 3726
$ws.Range("I136").Value = 4551.5
$ws.Range("K136").Value = 13654.5
$ws.Range("M136").Value = -11104.5
$ws.Range("H139").Value = 68845.92
$ws.Range("J139").Value = 88749.25
$ws.Range("L139").Value = 88749.25
$ws.Range("N139").Value = -99029.25
$ws.Range("H141").Value = 294245.62
$ws.Range("J141").Value = 294245.62
$ws.Range("L141").Value = 294245.62
$ws.Range("N141").Value = -304605.62

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
$ws.Range("H134").Value = 5113.4736
$ws.Range("I134").Value = 775.6429000000001
$ws.Range("J134").Value = 17259.4
$ws.Range("K134").Value = 2326.9287
$ws.Range("L134").Value = 51778.2
$ws.Range("M134").Value = 2743.0713
$ws.Range("N134").Value = -61918.2
$ws.Range("H136").Value = 6395.731
$ws.Range("I136").Value = 5254.6875
$ws.Range("J136").Value = 8221.4
$ws.Range("K136").Value = 15764.0625
$ws.Range("L136").Value = 24664.2
$ws.Range("M136").Value = -10664.0625
$ws.Range("N136").Value = -34864.2
$ws.Range("H138").Value = 861
$ws.Range("I138").Value = 861
$ws.Range("K138").Value = 2583
$ws.Range("M138").Value = 2557

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5899.579
$ws.Range("I80").Value = 3839
$ws.Range("J80").Value = 8189.1113
$ws.Range("K80").Value = 3839
$ws.Range("L80").Value = 8189.1113
$ws.Range("M80").Value = -2841
$ws.Range("N80").Value = -10185.1113
$ws.Range("H83").Value = 5899.579
$ws.Range("I83").Value = 3839
$ws.Range("J83").Value = 8189.1113
$ws.Range("K83").Value = 19195
$ws.Range("L83").Value = 40945.5565
$ws.Range("M83").Value = -14203
$ws.Range("N83").Value = -50929.5565
$ws.Range("H100").Value = 65000
$ws.Range("J100").Value = 65000
$ws.Range("L100").Value = 65000
$ws.Range("N100").Value = -67164
$ws.Range("H103").Value = 42333
$ws.Range("J103").Value = 42333
$ws.Range("L103").Value = 42333
$ws.Range("N103").Value = -44677
$ws.Range("H107").Value = 1094.8334
$ws.Range("I107").Value = 1119.5
$ws.Range("J107").Value = 1045.5
$ws.Range("K107").Value = 1119.5
$ws.Range("L107").Value = 1045.5
$ws.Range("M107").Value = 800.5
$ws.Range("N107").Value = -4885.5
$ws.Range("H132").Value = 4200.915
$ws.Range("I132").Value = 3892.575
$ws.Range("K132").Value = 11677.725
$ws.Range("M132").Value = -9147.724999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 15665.667
$ws.Range("J2").Value = 15665.667
$ws.Range("L2").Value = 15665.667
$ws.Range("N2").Value = -15889.667
$ws.Range("H46").Value = 6499.6665
$ws.Range("I46").Value = 4749.5
$ws.Range("K46").Value = 4749.5
$ws.Range("M46").Value = -4561.5
$ws.Range("H117").Value = 95696
$ws.Range("J117").Value = 95696
$ws.Range("L117").Value = 95696
$ws.Range("N117").Value = -104874
$ws.Range("H119").Value = 2433139.8
$ws.Range("J119").Value = 1649709.5
$ws.Range("L119").Value = 1649709.5
$ws.Range("N119").Value = -1659385.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 52833.332
$ws.Range("J98").Value = 52833.332
$ws.Range("L98").Value = 52833.332
$ws.Range("N98").Value = -58823.332
$ws.Range("H107").Value = 10258.8
$ws.Range("J107").Value = 10919.8
$ws.Range("L107").Value = 32759.4
$ws.Range("N107").Value = -36599.39999999999
$ws.Range("H132").Value = 6934
$ws.Range("I132").Value = 7117.5
$ws.Range("K132").Value = 21352.5
$ws.Range("M132").Value = -18822.5
$ws.Range("H136").Value = 3649.3547
$ws.Range("I136").Value = 3543.5
$ws.Range("K136").Value = 10630.5
$ws.Range("M136").Value = -8080.5
